# Command Strings List - update arm command strings, fix rover/drivetrain
# range labels, and renumber the Science section.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Rover Systems / Drive Train range labels (text corrections)
# ---------------------------------------------------------------------
$ws.Range("A4").Value = "Rover Systems (8000-8999)"
$ws.Range("A9").Value = "Drive Train (9000-9999)"
$ws.Range("D9").Value = "D"

# ---------------------------------------------------------------------
# 2. Arm section (rows 14-19): rewrite the DOF command table to match
#    the current (unfinished) firmware code.
# ---------------------------------------------------------------------
$ws.Range("A14").Value = "Arm (10000-10999)"

$ws.Range("B14").Value = "Shoulder Rotation"
$ws.Range("C14").Value = 10001
$ws.Range("D14").Value = "A"
$ws.Range("E14").Value = "R"
$ws.Range("F14").Value = "speed"
$ws.Range("G14").Value = "0 is stop, <0 CCW, >0 CW"
$ws.Range("H14").Value = "AR30"

$ws.Range("B15").Value = "Shoulder Angle"
$ws.Range("C15").Value = 10001
$ws.Range("D15").Value = "A"
$ws.Range("E15").Value = "L"
$ws.Range("F15").Value = "speed"
$ws.Range("G15").Value = "0 is stop, <0 down, >0 up"
$ws.Range("H15").Value = "AL30"

$ws.Range("B16").Value = "Elbow"
$ws.Range("C16").Value = 10001
$ws.Range("D16").Value = "A"
$ws.Range("E16").Value = "E"
$ws.Range("F16").Value = "speed"
$ws.Range("G16").Value = "0 is stop, <0 down, >0 up"
$ws.Range("H16").Value = "AE30"

$ws.Range("B17").Value = "Claw Pitch"
$ws.Range("C17").Value = 10001
$ws.Range("D17").Value = "C"
$ws.Range("E17").Value = "P"
$ws.Range("F17").Value = "speed"
$ws.Range("G17").Value = "0 is stop, <0 down, >0 up"
$ws.Range("H17").Value = "CP30"

$ws.Range("B18").Value = "Claw Rotation"
$ws.Range("C18").Value = 10001
$ws.Range("D18").Value = "C"
$ws.Range("E18").Value = "R"
$ws.Range("F18").Value = "speed"
$ws.Range("G18").Value = "0 is stop, <0 CCW, >0 CW"
$ws.Range("H18").Value = "CR30"

$ws.Range("B19").Value = "Claw Actuator"
$ws.Range("C19").Value = 10001
$ws.Range("D19").Value = "C"
$ws.Range("E19").Value = "C"
$ws.Range("F19").Value = "speed"
$ws.Range("G19").Value = "0 is stop, <0 retract, >0 extend"
$ws.Range("H19").Value = "CC30"

# rows 15-17 used to be tall (wrapped explanation text); the new text is
# short, so collapse them back to the normal row height used elsewhere
# in this block.
$ws.Rows.Item(14).RowHeight = 18.75
$ws.Rows.Item(15).RowHeight = 18.75
$ws.Rows.Item(16).RowHeight = 18.75
$ws.Rows.Item(17).RowHeight = 18.75
$ws.Rows.Item(18).RowHeight = 18.75
$ws.Rows.Item(19).RowHeight = 18.75

# ---------------------------------------------------------------------
# 3. Insert a blank spacer row before the Science section (old row 21),
#    matching the blank row above it (row 20).
# ---------------------------------------------------------------------
$ws.Rows.Item(21).Insert()
$ws.Rows.Item(21).RowHeight = 18.75

# ---------------------------------------------------------------------
# 4. Remove the now-redundant blank spacer row that used to sit right
#    after the Science section (old row 25, now shifted to row 26).
# ---------------------------------------------------------------------
$ws.Rows.Item(26).Delete()

# ---------------------------------------------------------------------
# 5. Re-add a trailing blank spacer row at the bottom of the sheet so
#    the table keeps the same amount of breathing room as before.
# ---------------------------------------------------------------------
$ws.Range("D29").HorizontalAlignment = -4108
$ws.Range("E29").HorizontalAlignment = -4108

# ---------------------------------------------------------------------
# 6. Window / selection cosmetics to match the saved view state.
# ---------------------------------------------------------------------
$ws.Application.ActiveWindow.ScrollRow = 10
$ws.Range("B19").Select()
